$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Row 12 (Testcase 11 - Application quit): add wait(5) before validate3 and
# change the minimized app check from the generic "homescreen" string to the
# concrete enterprise browser package name.
$ws.Range("G12").Value = "wait(3);`nvalidate1;`nlink_Click(PB2_test_link);`nwait(2);`nvalidate2;`nlink_Click(applicationpb_test_link);`nwait(5);`nvalidate3;"
$ws.Range("H12").Value = "validate1`n{`nvalidate_PageTitle=Pocket Browser Tests`n};`nvalidate2`n{`nvalidate_PageTitle=PB2 Tests`n};`nvalidate3`n{`nvalidate_AppMinimized=com.symbol.enterprisebrowser`n};"
$ws.Rows.Item(12).RowHeight = 147

# Row 15 (Testcase 14 - ScreenOrientation): extend the steps to go back and
# re-check normal orientation.
$ws.Range("G15").Value = "wait(3);`nvalidate1;`nlink_Click(PB2_test_link);`nwait(2);`nvalidate2;`nlink_Click(screenorientationpb_test_link);`nvalidate3;`npress_Key(Back);`nwait(2);`nlink_Click(screenorientationpbnormal_test_link);`nwait(5);"
$ws.Rows.Item(15).RowHeight = 166.5

# Row 20 (Testcase 19 - Quit button): fix the minimized app validation value.
$ws.Range("H20").Value = "validate1`n{`nvalidate_PageTitle=Pocket Browser Tests`n};`nvalidate2`n{`nvalidate_PageTitle=PB2 Tests`n};`nvalidate3`n{`nvalidate_isIconDisplayed=QuitButton_xpath,true`n};`nvalidate4`n{`nvalidate_isIconDisplayed=QuitButton_xpath,false`n};`nvalidate5`n{`nvalidate_AppMinimized=com.symbol.enterprisebrowser`n};"
$ws.Rows.Item(20).RowHeight = 294

# Row 22 (Testcase 21 - Minimize button): fix the minimized app validation value.
$ws.Range("H22").Value = "validate1`n{`nvalidate_PageTitle=Pocket Browser Tests`n};`nvalidate2`n{`nvalidate_PageTitle=PB2 Tests`n};`nvalidate3`n{`nvalidate_isIconDisplayed=MinimizeButton_xpath,true`n};`nvalidate4`n{`nvalidate_AppMinimized=com.symbol.enterprisebrowser`n};`nvalidate5`n{`nvalidate_isIconDisplayed=MinimizeButton_xpath,false`n};"
$ws.Rows.Item(22).RowHeight = 294
